$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")

# Row 33
$ws.Range("H33").Value = 489.8421
$ws.Range("I33").Value = 571.13336
$ws.Range("J33").Value = 185
$ws.Range("K33").Value = 571.13336
$ws.Range("L33").Value = 185
$ws.Range("M33").Value = -342.13336

# Row 51
$ws.Range("H51").Value = 4071.4167
$ws.Range("I51").Value = 6521.6
$ws.Range("J51").Value = 2321.2856
$ws.Range("K51").Value = 6521.6
$ws.Range("L51").Value = 2321.2856
$ws.Range("M51").Value = -6037.6
$ws.Range("N51").Value = -3289.2856

# Row 76
$ws.Range("H76").Value = 44536.75
$ws.Range("I76").Value = 44536.75
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 44536.75
$ws.Range("L76").Value = 0
$ws.Range("M76").Value = -44221.75
$ws.Range("N76").ClearContents()

# Row 79
$ws.Range("H79").Value = 44536.75
$ws.Range("I79").Value = 44536.75
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 44536.75
$ws.Range("L79").Value = 0
$ws.Range("M79").Value = -43444.75
$ws.Range("N79").ClearContents()

# Row 86
$ws.Range("H86").Value = 144258.86
$ws.Range("I86").Value = 167818.67
$ws.Range("J86").Value = 2900
$ws.Range("K86").Value = 167818.67
$ws.Range("L86").Value = 2900
$ws.Range("M86").Value = -166695.67
$ws.Range("N86").Value = -5146

# Row 89
$ws.Range("H89").Value = 144258.86
$ws.Range("I89").Value = 167818.67
$ws.Range("J89").Value = 2900
$ws.Range("K89").Value = 839093.3500000001
$ws.Range("L89").Value = 14500
$ws.Range("M89").Value = -833477.3500000001
$ws.Range("N89").Value = -25732

# Row 106
$ws.Range("H106").Value = 71573430
$ws.Range("I106").Value = 169001.5
$ws.Range("J106").Value = 500000000
$ws.Range("K106").Value = 169001.5
$ws.Range("L106").Value = 500000000
$ws.Range("M106").Value = -168370.5

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")

# Row 76
$ws.Range("H76").Value = 13000
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 13000
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 13000
$ws.Range("N76").Value = -13676

# Row 79
$ws.Range("H79").Value = 13000
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 13000
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 13000
$ws.Range("N79").Value = -15340

# Row 88
$ws.Range("H88").Value = 629614.4399999999
$ws.Range("I88").Value = 1003383.4
$ws.Range("J88").Value = 6666.1665
$ws.Range("K88").Value = 1003383.4
$ws.Range("L88").Value = 6666.1665
$ws.Range("M88").Value = -1002977.4
$ws.Range("N88").Value = -7478.1665

# Row 91
$ws.Range("H91").Value = 629614.4399999999
$ws.Range("I91").Value = 1003383.4
$ws.Range("J91").Value = 6666.1665
$ws.Range("K91").Value = 1003383.4
$ws.Range("L91").Value = 6666.1665
$ws.Range("M91").Value = -1001979.4
$ws.Range("N91").Value = -9474.166499999999

# Row 122
$ws.Range("H122").Value = 2112.3635
$ws.Range("I122").Value = 1212.2
$ws.Range("J122").Value = 2862.5
$ws.Range("K122").Value = 3636.6
$ws.Range("L122").Value = 8587.5
$ws.Range("M122").Value = -1186.6

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")

# Row 25
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("M25").ClearContents()

# Row 86
$ws.Range("H86").Value = 1593.9354
$ws.Range("I86").Value = 1500.0667
$ws.Range("J86").Value = 1681.9375
$ws.Range("K86").Value = 1500.0667
$ws.Range("L86").Value = 1681.9375
$ws.Range("M86").Value = -377.0667000000001
$ws.Range("N86").Value = -3927.9375

# Row 89
$ws.Range("H89").Value = 1593.9354
$ws.Range("I89").Value = 1500.0667
$ws.Range("J89").Value = 1681.9375
$ws.Range("K89").Value = 7500.333500000001
$ws.Range("L89").Value = 8409.6875
$ws.Range("M89").Value = -1884.333500000001
$ws.Range("N89").Value = -19641.6875

# Row 99
$ws.Range("H99").Value = 820
$ws.Range("I99").Value = 900
$ws.Range("J99").Value = 500
$ws.Range("K99").Value = 900
$ws.Range("L99").Value = 500
$ws.Range("M99").Value = 598
$ws.Range("N99").Value = -3496

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")

# Row 106
$ws.Range("H106").Value = 0
$ws.Range("I106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("K106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")

# Row 55
$ws.Range("H55").Value = 3090
$ws.Range("I55").Value = 1500
$ws.Range("J55").Value = 3266.6667
$ws.Range("K55").Value = 4500
$ws.Range("L55").Value = 9800.000100000001
$ws.Range("M55").Value = -4323
$ws.Range("N55").Value = -10154.0001

# Row 64
$ws.Range("H64").Value = 1660.9333
$ws.Range("I64").Value = 700
$ws.Range("J64").Value = 1808.7693
$ws.Range("K64").Value = 2100
$ws.Range("L64").Value = 5426.3079
$ws.Range("M64").Value = -1830
$ws.Range("N64").Value = -5966.3079

# Row 67
$ws.Range("H67").Value = 1660.9333
$ws.Range("I67").Value = 700
$ws.Range("J67").Value = 1808.7693
$ws.Range("K67").Value = 2100
$ws.Range("L67").Value = 5426.3079
$ws.Range("M67").Value = -1164
$ws.Range("N67").Value = -7298.3079

# Row 69
$ws.Range("H69").Value = 1300
$ws.Range("I69").Value = 300
$ws.Range("J69").Value = 1500
$ws.Range("K69").Value = 900
$ws.Range("L69").Value = 4500
$ws.Range("M69").Value = -89
$ws.Range("N69").Value = -6122

# Row 70
$ws.Range("H70").Value = 6249.75
$ws.Range("I70").Value = 4999
$ws.Range("J70").Value = 6666.6665
$ws.Range("K70").Value = 14997
$ws.Range("L70").Value = 19999.9995
$ws.Range("M70").Value = -14682

# Row 72
$ws.Range("H72").Value = 1300
$ws.Range("I72").Value = 300
$ws.Range("J72").Value = 1500
$ws.Range("K72").Value = 2700
$ws.Range("L72").Value = 13500
$ws.Range("M72").Value = 1356
$ws.Range("N72").Value = -21612

# Row 73
$ws.Range("H73").Value = 6249.75
$ws.Range("I73").Value = 4999
$ws.Range("J73").Value = 6666.6665
$ws.Range("K73").Value = 14997
$ws.Range("L73").Value = 19999.9995
$ws.Range("M73").Value = -13905

# Row 74
$ws.Range("H74").Value = 1990
$ws.Range("I74").Value = 1990
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 5970
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -4909
$ws.Range("N74").ClearContents()

# Row 75
$ws.Range("H75").Value = 1999.5
$ws.Range("I75").Value = 1999
$ws.Range("J75").Value = 2000
$ws.Range("K75").Value = 5997
$ws.Range("L75").Value = 6000
$ws.Range("M75").Value = -4999
$ws.Range("N75").Value = -7996

# Row 76
$ws.Range("H76").Value = 16632
$ws.Range("I76").Value = 2013
$ws.Range("J76").Value = 21505
$ws.Range("K76").Value = 6039
$ws.Range("L76").Value = 64515
$ws.Range("M76").Value = -5656
$ws.Range("N76").Value = -65281

# Row 77
$ws.Range("H77").Value = 1990
$ws.Range("I77").Value = 1990
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 17910
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -12606
$ws.Range("N77").ClearContents()

# Row 78
$ws.Range("H78").Value = 1999.5
$ws.Range("I78").Value = 1999
$ws.Range("J78").Value = 2000
$ws.Range("K78").Value = 17991
$ws.Range("L78").Value = 18000
$ws.Range("M78").Value = -12999
$ws.Range("N78").Value = -27984

# Row 79
$ws.Range("H79").Value = 16632
$ws.Range("I79").Value = 2013
$ws.Range("J79").Value = 21505
$ws.Range("K79").Value = 6039
$ws.Range("L79").Value = 64515
$ws.Range("M79").Value = -4713
$ws.Range("N79").Value = -67167

# Row 131
$ws.Range("H131").Value = 1434291.9
$ws.Range("I131").Value = 13088.75
$ws.Range("J131").Value = 2381760.8
$ws.Range("K131").Value = 39266.25
$ws.Range("L131").Value = 7145282.399999999
$ws.Range("M131").Value = -34226.25
$ws.Range("N131").Value = -7155362.399999999

# Row 132
$ws.Range("H132").Value = 896.125
$ws.Range("I132").Value = 863.8
$ws.Range("J132").Value = 950
$ws.Range("K132").Value = 7774.2
$ws.Range("L132").Value = 8550
$ws.Range("M132").Value = -5244.2
$ws.Range("N132").Value = -13610

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")

# Row 122
$ws.Range("H122").Value = 3618.1052
$ws.Range("I122").Value = 4468
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 13404
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -10954
$ws.Range("N122").Value = -13900

# Row 132
$ws.Range("H132").Value = 6548.375
$ws.Range("I132").Value = 8332.23
$ws.Range("J132").Value = 3235.5
$ws.Range("K132").Value = 24996.69
$ws.Range("L132").Value = 9706.5
$ws.Range("M132").Value = -22466.69
$ws.Range("N132").Value = -14766.5
